$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates, keyed by row number.
# D = $null means the Price column is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "61.645.76"; E = "  -2.34%  " },
    @{ Row = 3;  D = "2.898.67";  E = "  -1.97%  " },
    @{ Row = 4;  D = $null;       E = "  +0.04%  " },
    @{ Row = 5;  D = "568.67";    E = "  -4.46%  " },
    @{ Row = 6;  D = "143.48";    E = "  -3.90%  " },
    @{ Row = 7;  D = $null;       E = "  -0.01%  " },
    @{ Row = 8;  D = $null;       E = "  -0.54%  " },
    @{ Row = 9;  D = "2.896.09";  E = "  -2.03%  " },
    @{ Row = 10; D = $null;       E = "  -2.00%  " },
    @{ Row = 11; D = $null;       E = "  -2.59%  " },
    @{ Row = 12; D = $null;       E = "  -2.56%  " },
    @{ Row = 13; D = $null;       E = "  -1.75%  " },
    @{ Row = 14; D = "31.91";     E = "  -2.96%  " },
    @{ Row = 15; D = $null;       E = "  -0.55%  " },
    @{ Row = 16; D = "3.377.43";  E = "  -2.05%  " },
    @{ Row = 17; D = "61.623.79"; E = "  -2.29%  " },
    @{ Row = 18; D = "6.55";      E = "  -2.32%  " },
    @{ Row = 19; D = "2.908.64";  E = "  -1.67%  " },
    @{ Row = 20; D = "433.39";    E = "  -2.11%  " },
    @{ Row = 21; D = "13.09";     E = "  -3.28%  " },
    @{ Row = 22; D = "0.658";     E = "  -1.80%  " },
    @{ Row = 23; D = $null;       E = "  -2.84%  " },
    @{ Row = 24; D = "79.48";     E = "  -1.92%  " },
    @{ Row = 25; D = "12.00";     E = "  +1.68%  " },
    @{ Row = 26; D = $null;       E = "  -0.03%  " },
    @{ Row = 27; D = "9.95";      E = "  -11.20%  " },
    @{ Row = 28; D = $null;       E = "  -6.13%  " },
    @{ Row = 29; D = $null;       E = "  +4.32%  " },
    @{ Row = 30; D = "7.00";      E = "  -4.68%  " },
    @{ Row = 31; D = $null;       E = "  -4.17%  " },
    @{ Row = 32; D = $null;       E = "  -8.07%  " },
    @{ Row = 33; D = $null;       E = "  +0.07%  " },
    @{ Row = 34; D = $null;       E = "  -2.08%  " },
    @{ Row = 35; D = "25.58";     E = "  -3.51%  " },
    @{ Row = 36; D = "0.960";     E = "  -3.28%  " },
    @{ Row = 37; D = $null;       E = "  -3.99%  " },
    @{ Row = 38; D = "48.87";     E = "  -1.89%  " },
    @{ Row = 39; D = $null;       E = "  -5.78%  " },
    @{ Row = 40; D = $null;       E = "  -9.07%  " },
    @{ Row = 41; D = $null;       E = "  -3.72%  " },
    @{ Row = 42; D = $null;       E = "  -3.08%  " },
    @{ Row = 43; D = "39.38";     E = "  +0.36%  " },
    @{ Row = 44; D = $null;       E = "  -5.27%  " },
    @{ Row = 45; D = "2.690.24";  E = "  -0.44%  " },
    @{ Row = 46; D = "133.12";    E = "  -1.74%  " },
    @{ Row = 47; D = "0.0335";    E = "  -0.95%  " },
    @{ Row = 48; D = "347.21";    E = "  -3.83%  " },
    @{ Row = 50; D = $null;       E = "  -1.74%  " },
    @{ Row = 51; D = $null;       E = "  -5.49%  " }
)

# Rows whose new Price text would otherwise be auto-parsed into a numeric
# value by Excel's smart-typing (single "." as decimal separator). Force
# the cell to stay text by temporarily marking it as Text ("@") before
# assigning the value, then restore its (unstyled) appearance.
$forceTextRows = @(5, 6, 14, 18, 20, 21, 22, 24, 25, 27, 30, 35, 36, 38, 43, 46, 47, 48)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($forceTextRows -contains $r) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}

foreach ($r in $forceTextRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
